# Corrected Scaling and Sampling
# Both scaling and oversampling (when applicable) are now applied after the train/val split.
# Re-running feature selection (RFE/RFECV/Logistics/RandomForest/LightGBM/Lasso/Ridge/Elastic)
# with the corrected pipeline produced slightly different per-model vote counts/rankings for some
# mid-table features (the top and bottom of each ranking are unaffected). Apply those updated
# results to the two results sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "final_fail" ---
$ws1 = $wb.Worksheets.Item("final_fail")

# Row 16
$ws1.Range("A16").Value = "Start of Session 2 (%)"
$ws1.Range("B16").Value = $true
$ws1.Range("D16").Value = $false

# Row 17
$ws1.Range("A17").Value = "Number of clicks"
$ws1.Range("B17").Value = $false
$ws1.Range("D17").Value = $true
$ws1.Range("E17").Value = $true
$ws1.Range("J17").Value = 5

# Row 18
$ws1.Range("A18").Value = "Start of Session 7 (%)"

# Row 19
$ws1.Range("A19").Value = "Clicks on course"
$ws1.Range("B19").Value = $true
$ws1.Range("E19").Value = $false

# Row 20
$ws1.Range("A20").Value = "Days with no interaction (%)"
$ws1.Range("B20").Value = $false
$ws1.Range("E20").Value = $true

# Row 22
$ws1.Range("A22").Value = "Start of Session 6 (%)"
$ws1.Range("E22").Value = $false
$ws1.Range("J22").Value = 3

# Row 23
$ws1.Range("A23").Value = "Clicks per day"
$ws1.Range("E23").Value = $true
$ws1.Range("F23").Value = $false

# Row 24
$ws1.Range("A24").Value = "Clicks on folder"
$ws1.Range("D24").Value = $true
$ws1.Range("E24").Value = $false

# Row 25
$ws1.Range("A25").Value = "Assignments submitted"

# Row 26
$ws1.Range("A26").Value = "Start of Session 4 (%)"
$ws1.Range("D26").Value = $false
$ws1.Range("F26").Value = $true

# Row 28
$ws1.Range("A28").Value = "Clicks on forum"
$ws1.Range("C28").Value = $true
$ws1.Range("D28").Value = $false

# Row 33
$ws1.Range("A33").Value = "Start of Session 8 (%)"
$ws1.Range("C33").Value = $true
$ws1.Range("J33").Value = 2

# Row 34
$ws1.Range("A34").Value = "Start of Session 5 (%)"
$ws1.Range("C34").Value = $true
$ws1.Range("J34").Value = 2

# Row 35
$ws1.Range("A35").Value = "Links viewed"
$ws1.Range("D35").Value = $true
$ws1.Range("J35").Value = 2

# Row 36
$ws1.Range("A36").Value = "Start of Session 9 (%)"

# --- Sheet "final_gifted" ---
$ws2 = $wb.Worksheets.Item("final_gifted")

# Row 3
$ws2.Range("A3").Value = "Total time online (min)"

# Row 7
$ws2.Range("A7").Value = "Clicks per session"

# Row 8
$ws2.Range("A8").Value = "Clicks (% of course total)"
$ws2.Range("D8").Value = $true
$ws2.Range("J8").Value = 6

# Row 9
$ws2.Range("A9").Value = "Largest period of inactivity (h)"
$ws2.Range("C9").Value = $true
$ws2.Range("D9").Value = $false

# Row 11
$ws2.Range("A11").Value = "Start of Session 7 (%)"
$ws2.Range("E11").Value = $false
$ws2.Range("J11").Value = 4

# Row 13
$ws2.Range("A13").Value = "Clicks per day"
$ws2.Range("C13").Value = $true
$ws2.Range("F13").Value = $false

# Row 14
$ws2.Range("A14").Value = "Submissions (% of course total)"
$ws2.Range("E14").Value = $false
$ws2.Range("F14").Value = $true

# Row 15
$ws2.Range("A15").Value = "Start of Session 4 (%)"
$ws2.Range("E15").Value = $true
$ws2.Range("J15").Value = 4

# Row 16
$ws2.Range("A16").Value = "Start of Session 3 (%)"
$ws2.Range("B16").Value = $true
$ws2.Range("J16").Value = 4

# Row 17
$ws2.Range("A17").Value = "Assignments viewed"
$ws2.Range("E17").Value = $true
$ws2.Range("F17").Value = $false

# Row 18
$ws2.Range("A18").Value = "Days with no interaction (%)"
$ws2.Range("B18").Value = $false
$ws2.Range("F18").Value = $true

# Row 19
$ws2.Range("A19").Value = "Number of days"
$ws2.Range("B19").Value = $true
$ws2.Range("E19").Value = $false

# Row 22
$ws2.Range("A22").Value = "Number of clicks"
$ws2.Range("B22").Value = $false
$ws2.Range("E22").Value = $true

# Row 27
$ws2.Range("A27").Value = "Start of Session 6 (%)"
